# Automatische test-sync: 2025-07-22 12:44:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 15 to the Logs sheet
$logs.Range("A15").Value = "Wat is jullie privacybeleid?"
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value = "Testmail #15: Wat is jullie privacybeleid?"
$logs.Range("D15").Value = "Overig"
$logs.Range("E15").Value = "Beste afzender,
Dank u voor uw vraag over ons privacybeleid. Wij hechten veel waarde aan de bescherming van persoonlijke gegevens van onze klanten en volgen strikte richtlijnen om deze te waarborgen. Ons privacybeleid is te vinden op onze website onder [link naar privacybeleid]. Hier vindt u gedetailleerde informatie over hoe wij omgaan met persoonlijke gegevens, welke gegevens wij verzamelen, hoe wij deze gebruiken en welke maatregelen wij treffen om ze te beschermen.
Mocht u nog verdere vragen hebben over ons privacybeleid, dan helpen wij graag verder.
Met vriendelijke groet,
[Naam] 
E-mailassistent"
$logs.Range("F15").Value = "2025-07-22 12:43:56"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Ja"
$logs.Range("J15").Value = "Ja"

# Writing multi-line text into a brand-new row triggers a simulated
# Excel "fit row to wrapped text" autosize (ht/customHeight on <row>).
# Re-running AutoFit afterwards clears that custom-height flag again,
# matching the openpyxl-authored source (no explicit row height).
$logs.Rows.Item(15).AutoFit()

# Extend the conditional-formatting ranges from row 14 to the new row 15
$logs.Range("D2:D14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))
$logs.Range("G2:G14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))
$logs.Range("H2:H14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H15"))
$logs.Range("I2:I14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I15"))
$logs.Range("J2:J14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J15"))

# Update the Dashboard summary: swap "Bestelling / Levering" and "Overig" rows,
# and bump the "Overig" count for the new privacybeleid entry
$dashboard.Range("A4").Value = "Overig"
$dashboard.Range("B4").Value = 2
$dashboard.Range("A6").Value = "Bestelling / Levering"
$dashboard.Range("B6").Value = 1
